# Commit: "Updated PowerPoint record description and code formatting."
#
# Slide 5 ("Records/Immutability") -> Content Placeholder 2 -> last bullet:
#   "Reference type with value semantics"
# The tail of the sentence is reworded, splitting the single run into two:
#   "Reference type with " (kept as-is)
#   "value equality"       (new wording)

$p     = $ppt.ActivePresentation
$slide = $p.Slides.Item(5)
$shape = $slide.Shapes.Item("Content Placeholder 2")
$body  = $shape.TextFrame.TextRange

# Locate "value semantics" inside the bullet text and replace just that
# portion, leaving "Reference type with " untouched. Setting .Text on the
# found sub-range naturally splits the paragraph's run at that boundary.
$found = $body.Find("value semantics", 0)
$found.Text = "value equality"
